# Fix property_category mislabeling:
#  - 建物 (Building) sheet: rows 2-10, column I ("property_category") were
#    tagged "land" but should be "building".
#  - 汽車 (Car) sheet: row 2, column H ("property_category") was tagged
#    "land" but should be "car".

$wb = $excel.ActiveWorkbook

$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I10").Value = "building"

$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
